$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-ambiguous text values (won't be auto-converted to numbers) - direct assignment
# Ambiguous numeric-looking text values - use NumberFormat "@" + Style reset trick to force text

$ws.Range("D2").Value = "65.860.11"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "2.665.55"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.653"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.38%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -2.22%  "
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("D15").Value = "3.146.86"
$ws.Range("E15").Value = "  -0.33%  "
$ws.Range("D16").Value = "65.770.70"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "2.664.29"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "351.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("E21").Value = "  -1.12%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "576.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.55%  "
$ws.Range("E30").Value = "  -2.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.43%  "
$ws.Range("E35").Value = "  +1.76%  "
$ws.Range("E36").Value = "  -0.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.58"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "154.18"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "161.29"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.39%  "
$ws.Range("E43").Value = "  +1.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.644"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0257"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.19%  "
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("E50").Value = "  -6.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.815"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.37%  "
